$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39 is the grand total row for the San Diego area. Relabel it:
#   A39: "All other ports" -> "SAN DIEGO AREA TOTALS"
#   C39: "SAN DIEGO AREA TOTALS" -> "Totals"
$ws.Range("A39").Value = "SAN DIEGO AREA TOTALS"
$ws.Range("C39").Value = "Totals"

# Column A now needs to be widened to fit the longer label (matches the
# existing best-fit width already used by column C, which held this same
# text before the edit).
$ws.Columns("A").ColumnWidth = 21.833333333333332

# Select the whole of column A, mirroring the user clicking on the column
# header after making the edit.
$ws.Columns("A").Select() | Out-Null
